$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update recalculated statistics in the existing sheets
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("descriptives")
$ws.Range("C2").Value = [double]"845"
$ws.Range("D2").Value = [double]"0.008095743408439994"
$ws.Range("E2").Value = [double]"0.05052704415099526"
$ws.Range("F2").Value = [double]"97.73316103332681"
$ws.Range("G2").Value = [double]"84.23631744806789"
$ws.Range("H2").Value = [double]"13.49684358525892"
$ws.Range("C3").Value = [double]"438"
$ws.Range("D3").Value = [double]"0.002763808288530445"
$ws.Range("E3").Value = [double]"0.05318914704310809"
$ws.Range("F3").Value = [double]"95.02046485971272"
$ws.Range("G3").Value = [double]"90.32690851755481"
$ws.Range("H3").Value = [double]"4.693556342157903"

$ws = $wb.Worksheets.Item("coefficients")
$ws.Range("D2").Value = [double]"-0.1597154571624137"
$ws.Range("E2").Value = [double]"2.269635605788456e-14"
$ws.Range("F2").Value = [double]"-7097821819688.762"
$ws.Range("G2").Value = [double]"8.969227300150751e-14"
$ws.Range("H2").Value = [double]"-0.1597154571626948"
$ws.Range("I2").Value = [double]"-0.1597154571621327"
$ws.Range("D3").Value = [double]"0.1213364039631816"
$ws.Range("E3").Value = [double]"0.01855766108550308"
$ws.Range("F3").Value = [double]"6.570719159327477"
$ws.Range("G3").Value = [double]"1.774098066587662e-08"
$ws.Range("H3").Value = [double]"0.08455854362394009"
$ws.Range("I3").Value = [double]"0.1577841082989081"
$ws.Range("J3").Value = [double]"55.94756359842736"
$ws.Range("D4").Value = [double]"0.003972252781116893"
$ws.Range("E4").Value = [double]"6.293566203749799e-05"
$ws.Range("F4").Value = [double]"63.11641993053626"
$ws.Range("G4").Value = [double]"0.01008559281622121"
$ws.Range("H4").Value = [double]"0.003172589622263816"
$ws.Range("I4").Value = [double]"0.004771910859719281"
$ws.Range("D5").Value = [double]"0.1229849656393078"
$ws.Range("E5").Value = [double]"0.06178650348428123"
$ws.Range("F5").Value = [double]"2.000610310205984"
$ws.Range("G5").Value = [double]"0.2950895471579177"
$ws.Range("H5").Value = [double]"-0.5793350450826276"
$ws.Range("I5").Value = [double]"0.7204993810036464"
$ws.Range("J5").Value = [double]"0.9999999999999982"
$ws.Range("D6").Value = [double]"0.0480027512630412"
$ws.Range("E6").Value = [double]"0.236169282017018"
$ws.Range("F6").Value = [double]"0.2034120282344944"
$ws.Range("G6").Value = [double]"0.8722468408253656"
$ws.Range("H6").Value = [double]"-0.994566243336082"
$ws.Range("I6").Value = [double]"0.9955140641870832"
$ws.Range("D7").Value = [double]"0.03454035753874145"
$ws.Range("E7").Value = [double]"0.0205491837863112"
$ws.Range("F7").Value = [double]"1.681531670830845"
$ws.Range("G7").Value = [double]"0.1250473182328547"
$ws.Range("H7").Value = [double]"-0.01153354948722172"
$ws.Range("I7").Value = [double]"0.08046791089556368"
$ws.Range("J7").Value = [double]"9.538981409281005"
$ws.Range("D8").Value = [double]"0.04491997566670915"
$ws.Range("E8").Value = [double]"2.221481368528577e-15"
$ws.Range("F8").Value = [double]"20234347302518.88"
$ws.Range("G8").Value = [double]"3.146233297519468e-14"
$ws.Range("H8").Value = [double]"0.04491997566668098"
$ws.Range("I8").Value = [double]"0.04491997566673731"
$ws.Range("D9").Value = [double]"0.03775206423673266"
$ws.Range("E9").Value = [double]"0.01870525575479194"
$ws.Range("F9").Value = [double]"2.01921936072586"
$ws.Range("G9").Value = [double]"0.05102290451872862"
$ws.Range("H9").Value = [double]"-0.0001772413081629961"
$ws.Range("I9").Value = [double]"0.07557290314991361"
$ws.Range("J9").Value = [double]"35.69552461908427"
$ws.Range("D10").Value = [double]"0.1849994946729422"
$ws.Range("E10").Value = [double]"0.08723663062783917"
$ws.Range("F10").Value = [double]"2.145365376706126"
$ws.Range("G10").Value = [double]"0.2776801997575834"
$ws.Range("H10").Value = [double]"-0.7265080720437419"
$ws.Range("I10").Value = [double]"0.8605863923566883"
$ws.Range("J10").Value = [double]"1"
$ws.Range("D11").Value = [double]"0.1091262032705729"
$ws.Range("E11").Value = [double]"0.04309853415595803"
$ws.Range("F11").Value = [double]"2.542139879805543"
$ws.Range("G11").Value = [double]"0.05881625065882558"
$ws.Range("H11").Value = [double]"-0.006307984203805485"
$ws.Range("I11").Value = [double]"0.2216902961138924"
$ws.Range("J11").Value = [double]"4.359205059832483"

$ws = $wb.Worksheets.Item("pairwise")
$ws.Range("C2").Value = [double]"15.25148396653277"
$ws.Range("E2").Value = [double]"55.94756359842755"
$ws.Range("F2").Value = [double]"1.384407059621348e-21"
$ws.Range("G2").Value = [double]"2.076610589432022e-20"
$ws.Range("C3").Value = [double]"2622.789044953463"
$ws.Range("F3").Value = [double]"0.0002427262469866367"
$ws.Range("G3").Value = [double]"0.0009102234261998876"
$ws.Range("C4").Value = [double]"4.607889929064195"
$ws.Range("E4").Value = [double]"0.9999999999999984"
$ws.Range("F4").Value = [double]"0.1360489697985592"
$ws.Range("G4").Value = [double]"0.3401224244963979"
$ws.Range("C5").Value = [double]"0.88552737312287"
$ws.Range("F5").Value = [double]"0.53860256822631"
$ws.Range("G5").Value = [double]"0.807903852339465"
$ws.Range("C6").Value = [double]"9.521000767525148"
$ws.Range("E6").Value = [double]"9.538981409280998"
$ws.Range("F6").Value = [double]"3.527848647479622e-06"
$ws.Range("G6").Value = [double]"1.763924323739811e-05"
$ws.Range("C7").Value = [double]"6.356632264173754"
$ws.Range("E7").Value = [double]"55.96196981240531"
$ws.Range("F7").Value = [double]"3.983221008873925e-08"
$ws.Range("G7").Value = [double]"2.987415756655444e-07"
$ws.Range("C8").Value = [double]"0.02594097773135135"
$ws.Range("E8").Value = [double]"1.040428343873796"
$ws.Range("F8").Value = [double]"0.98336390185164"
$ws.Range("G8").Value = [double]"0.98336390185164"
$ws.Range("C9").Value = [double]"0.3119390382671925"
$ws.Range("E9").Value = [double]"1.074758992053866"
$ws.Range("F9").Value = [double]"0.8046708665762504"
$ws.Range("G9").Value = [double]"0.98336390185164"
$ws.Range("C10").Value = [double]"3.155928332306659"
$ws.Range("E10").Value = [double]"12.5240851518328"
$ws.Range("F10").Value = [double]"0.007899642622595907"
$ws.Range("G10").Value = [double]"0.02369892786778772"
$ws.Range("C11").Value = [double]"1.936318992284431"
$ws.Range("E11").Value = [double]"1.00000517278386"
$ws.Range("F11").Value = [double]"0.3034858176971905"
$ws.Range("G11").Value = [double]"0.5690359081822322"
$ws.Range("C12").Value = [double]"0.1865924181457022"
$ws.Range("E12").Value = [double]"1.000009488474191"
$ws.Range("F12").Value = [double]"0.8825618621400841"
$ws.Range("G12").Value = [double]"0.98336390185164"
$ws.Range("C13").Value = [double]"1.488219024285826"
$ws.Range("E13").Value = [double]"9.539346128281155"
$ws.Range("F13").Value = [double]"0.1689863212889058"
$ws.Range("G13").Value = [double]"0.3621135456190838"
$ws.Range("C14").Value = [double]"0.309567931083228"
$ws.Range("E14").Value = [double]"1.840522572537709"
$ws.Range("F14").Value = [double]"0.7883838335059444"
$ws.Range("G14").Value = [double]"0.98336390185164"
$ws.Range("C15").Value = [double]"1.367701705324735"
$ws.Range("E15").Value = [double]"1.28639194939222"
$ws.Range("F15").Value = [double]"0.362444462914999"
$ws.Range("G15").Value = [double]"0.6040741048583317"
$ws.Range("C16").Value = [double]"0.05688635245826114"
$ws.Range("E16").Value = [double]"1.547872764138074"
$ws.Range("F16").Value = [double]"0.9610947947159141"
$ws.Range("G16").Value = [double]"0.98336390185164"
$ws.Range("C17").Value = [double]"0.3838606145568483"
$ws.Range("E17").Value = [double]"35.69552461908432"
$ws.Range("F17").Value = [double]"0.7033588198889126"
$ws.Range("G17").Value = [double]"0.7033588198889126"
$ws.Range("C18").Value = [double]"1.630097590532166"
$ws.Range("E18").Value = [double]"1"
$ws.Range("F18").Value = [double]"0.3503051007771155"
$ws.Range("G18").Value = [double]"0.5254576511656732"
$ws.Range("C19").Value = [double]"1.499175741530596"
$ws.Range("E19").Value = [double]"4.359205059832481"
$ws.Range("F19").Value = [double]"0.2024454607297748"
$ws.Range("G19").Value = [double]"0.5254576511656732"
$ws.Range("C20").Value = [double]"1.674347639661599"
$ws.Range("E20").Value = [double]"1.039284206784907"
$ws.Range("F20").Value = [double]"0.3357635555885418"
$ws.Range("G20").Value = [double]"0.5254576511656732"
$ws.Range("C21").Value = [double]"1.528063534952901"
$ws.Range("E21").Value = [double]"5.256919591404019"
$ws.Range("F21").Value = [double]"0.1842160499475272"
$ws.Range("G21").Value = [double]"0.5254576511656732"
$ws.Range("C22").Value = [double]"0.7974327813315768"
$ws.Range("E22").Value = [double]"1.420021419896588"
$ws.Range("F22").Value = [double]"0.5366515439152819"
$ws.Range("G22").Value = [double]"0.6439818526983383"

# ---------------------------------------------------------------
# 2) Add the new "nr_studies" sheet at the end of the workbook
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ns = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ns.Name = "nr_studies"

$ns.Cells.Item(1,1).Value = "outcome"
$ns.Cells.Item(1,2).Value = "moderator_design"
$ns.Cells.Item(1,3).Value = "n_effect_sizes"
$ns.Cells.Item(1,4).Value = "k_studies"
$ns.Range("A1:D1").Font.Bold = $true
$ns.Range("A1:D1").HorizontalAlignment = -4108

$ns.Cells.Item(2,1).Value = "NS"
$ns.Cells.Item(2,2).Value = "Cross-sectional"
$ns.Cells.Item(2,3).Value = 732
$ns.Cells.Item(2,4).Value = 65

$ns.Cells.Item(3,1).Value = "NS"
$ns.Cells.Item(3,2).Value = "Experimental (non-randomized)"
$ns.Cells.Item(3,3).Value = 21
$ns.Cells.Item(3,4).Value = 2

$ns.Cells.Item(4,1).Value = "NS"
$ns.Cells.Item(4,2).Value = "Longitudinal"
$ns.Cells.Item(4,3).Value = 65
$ns.Cells.Item(4,4).Value = 12

$ns.Cells.Item(5,1).Value = "NS"
$ns.Cells.Item(5,2).Value = "Experimental (RCT)"
$ns.Cells.Item(5,3).Value = 13
$ns.Cells.Item(5,4).Value = 1

$ns.Cells.Item(6,1).Value = "NS"
$ns.Cells.Item(6,2).Value = "Cross-lagged"
$ns.Cells.Item(6,3).Value = 12
$ns.Cells.Item(6,4).Value = 1

$ns.Cells.Item(7,1).Value = "NT"
$ns.Cells.Item(7,2).Value = "Cross-sectional"
$ns.Cells.Item(7,3).Value = 395
$ns.Cells.Item(7,4).Value = 42

$ns.Cells.Item(8,1).Value = "NT"
$ns.Cells.Item(8,2).Value = "Experimental (non-randomized)"
$ns.Cells.Item(8,3).Value = 12
$ns.Cells.Item(8,4).Value = 2

$ns.Cells.Item(9,1).Value = "NT"
$ns.Cells.Item(9,2).Value = "Longitudinal"
$ns.Cells.Item(9,3).Value = 23
$ns.Cells.Item(9,4).Value = 6

$ns.Cells.Item(10,1).Value = "NT"
$ns.Cells.Item(10,2).Value = "Cross-lagged"
$ns.Cells.Item(10,3).Value = 8
$ns.Cells.Item(10,4).Value = 1

$ns.Cells.Item(11,1).Value = "NS"
$ns.Cells.Item(11,2).Value = "Experimental (non-randomized"
$ns.Cells.Item(11,3).Value = 2
$ns.Cells.Item(11,4).Value = 1
